$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '41.514.31'
$ws.Range("E2").Value = '  -2.34%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.485.81'
$ws.Range("E3").Value = '  -1.30%  '
$ws.Range("E4").Value = '  +0.21%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.65'
$ws.Range("E5").Value = '  +0.64%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '94.53'
$ws.Range("E6").Value = '  -4.33%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.549'
$ws.Range("E7").Value = '  -2.46%  '
$ws.Range("E8").Value = '  +0.15%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.500'
$ws.Range("E9").Value = '  -3.31%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '33.61'
$ws.Range("E10").Value = '  -4.50%  '
$ws.Range("E11").Value = '  -2.38%  '
$ws.Range("E12").Value = '  -0.16%  '
$ws.Range("E13").Value = '  -2.91%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.871.08'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.51'
$ws.Range("E15").Value = '  +1.52%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.475.96'
$ws.Range("E16").Value = '  -4.76%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.793'
$ws.Range("E17").Value = '  -1.81%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '41.473.74'
$ws.Range("E18").Value = '  -2.38%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.36'
$ws.Range("E19").Value = '  -3.64%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0924'
$ws.Range("E20").Value = '  -2.16%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.31'
$ws.Range("E21").Value = '  -5.19%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.82'
$ws.Range("E22").Value = '  -0.09%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '237.19'
$ws.Range("E23").Value = '  -1.43%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.76'
$ws.Range("E24").Value = '  -2.88%  '
$ws.Range("B25").Value = 'Dai'
$ws.Range("C25").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  -0.04%  '
$ws.Range("B26").Value = 'ImmutableX'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.90'
$ws.Range("E26").Value = '  -4.20%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.16'
$ws.Range("E27").Value = '  -4.07%  '
$ws.Range("E28").Value = '  -1.17%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.77'
$ws.Range("E29").Value = '  -2.54%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.67'
$ws.Range("E30").Value = '  -4.53%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '152.25'
$ws.Range("E31").Value = '  -2.62%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.50'
$ws.Range("E32").Value = '  -6.12%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.59'
$ws.Range("E33").Value = '  -3.19%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '18.12'
$ws.Range("E34").Value = '  +3.89%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0757'
$ws.Range("E35").Value = '  -3.48%  '
$ws.Range("B36").Value = 'ApeXProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.50'
$ws.Range("E36").Value = '  -10.68%  '
$ws.Range("B37").Value = 'LidoDAOToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.09'
$ws.Range("E37").Value = '  -1.29%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.88'
$ws.Range("E38").Value = '  -3.73%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.114'
$ws.Range("E39").Value = '  -1.77%  '
$ws.Range("E40").Value = '  -5.41%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.23'
$ws.Range("E41").Value = '  +1.52%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.01'
$ws.Range("E42").Value = '  +0.29%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '19.82'
$ws.Range("E43").Value = '  -8.14%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.992.32'
$ws.Range("E44").Value = '  -0.22%  '
$ws.Range("E45").Value = '  -2.77%  '
$ws.Range("E46").Value = '  -6.80%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.86'
$ws.Range("E47").Value = '  -2.39%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.735.13'
$ws.Range("E48").Value = '  -0.80%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '69.89'
$ws.Range("E49").Value = '  -1.96%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '97.11'
$ws.Range("E50").Value = '  -3.09%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.178'
$ws.Range("E51").Value = '  -5.33%  '
